# Hortaliza, Vega Central Mapocho de Santiago - Melón
# Insert 6 new daily-price rows (date 45008) right before the existing
# row 982 block, pushing all subsequent rows down by 6 (982-1064 -> 988-1070).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything currently at row 982 and below down by six rows.
$ws.Rows("982:987").Insert()

# Constant columns shared by every row in this market/category block.
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$catId     = 100112027
$categoria = "Melón"
$unidadCom = "`$/unidad"
$origen    = "Región de O'Higgins"
$kgOUnidad = 1
$clasif    = "Hortaliza"
$fecha     = 45008

# NOTE: named PowerShell parameters (-Name value) don't bind reliably in
# this COM host, so Set-Fila is called positionally.
function Set-Fila($Fila, $Variedad, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg) {
    $ws.Cells.Item($Fila, 1).Value  = $mercadoId
    $ws.Cells.Item($Fila, 2).Value  = $mercado
    $ws.Cells.Item($Fila, 3).Value  = $region
    $ws.Cells.Item($Fila, 4).Value  = $fecha
    $ws.Cells.Item($Fila, 5).Value  = $codreg
    $ws.Cells.Item($Fila, 6).Value  = $catId
    $ws.Cells.Item($Fila, 7).Value  = $categoria
    $ws.Cells.Item($Fila, 8).Value  = $Variedad
    $ws.Cells.Item($Fila, 9).Value  = $Calidad
    $ws.Cells.Item($Fila, 10).Value = $Volumen
    $ws.Cells.Item($Fila, 11).Value = $PrecioMin
    $ws.Cells.Item($Fila, 12).Value = $PrecioMax
    $ws.Cells.Item($Fila, 13).Value = $PrecioProm
    $ws.Cells.Item($Fila, 14).Value = $unidadCom
    $ws.Cells.Item($Fila, 15).Value = $origen
    $ws.Cells.Item($Fila, 16).Value = $PrecioKg
    $ws.Cells.Item($Fila, 17).Value = $kgOUnidad
    $ws.Cells.Item($Fila, 18).Value = $clasif
}

Set-Fila 982 "Calameño" "Primera" 430 900 900 900 900
Set-Fila 983 "Calameño" "Segunda" 340 700 700 700 700
Set-Fila 984 "Calameño" "Tercera" 160 500 500 500 500
Set-Fila 985 "Tuna"     "Primera" 520 900 900 900 900
Set-Fila 986 "Tuna"     "Segunda" 430 700 700 700 700
Set-Fila 987 "Tuna"     "Tercera" 250 500 500 500 500
